$d = $word.ActiveDocument

# Anchor: the paragraph right before the final empty bookmark paragraph
# (the bullet ending "...Known as random, or monte carlo algorithm. ")
$anchorIndex = $d.Paragraphs.Count - 1
$anchorRange = $d.Paragraphs.Item($anchorIndex).Range
$anchorRange.Collapse(0)

function Insert-ParagraphXml($afterRange, $xmlFragment) {
    $afterRange.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $newRange = $newPara.Range
    $newRange.InsertXML($xmlFragment)
    $newRange.Collapse(0)
    return $newRange
}

$xml0 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Addition in binary is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(^), and multiplication in binary is logical and (&amp;)</w:t></w:r></w:p>'
$xml1 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>One-time Pad</w:t></w:r><w:r><w:t>. plaintext XOR key = ciphertext, and ciphertext XOR key = plaintext.</w:t></w:r></w:p>'
$xml2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>If you XOR two one-time pad ciphertexts, you can recover the key by</w:t></w:r></w:p>'
$xml3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Blum-Blum-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shub</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> good, but computationally expensive. LFSR (Linear Feedback Shift Register) is less secure, less random, but much faster.</w:t></w:r></w:p>'
$xml4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>A flipflop (FF) is a memory device that holds one bit</w:t></w:r></w:p>'
$xml5 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>a | b means a divides b, which is true if and only if b % a = 0 (a goes into b)</w:t></w:r></w:p>'
$xml6 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Prime number theorem (density): number of primes &lt; x is approximately x/ln(x)</w:t></w:r><w:r><w:t xml:space="preserve"> ...generally</w:t></w:r></w:p>'
$xml7 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Fundamental theorem of arithmetic</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Every positive integer can be uniquely factored into prime numbers: no two prime factorizations are alike...they act as fingerprints</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r></w:p>'
$xml8 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">To find a-inverse, mod b: ax + by = 1. by mod b - 0; 1 mod b = 1; </w:t></w:r></w:p>'
$xml9 = '<w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/><w:t>ax + by (mod b) is ax = 1(mod b)</w:t></w:r></w:p>'
$xml10 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>DES is a symmetric (private) block cipher that encrypts 64-bit (8-byte, 8-character) blocks</w:t></w:r></w:p>'
$xml11 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Can obtain non-linearity not only by raising to power, but also by writing functions/operations so that distributive property doesn''t hold (DES, for example)</w:t></w:r></w:p>'

$cur = $anchorRange
$cur = Insert-ParagraphXml $cur $xml0
$cur = Insert-ParagraphXml $cur $xml1
$cur = Insert-ParagraphXml $cur $xml2
$cur = Insert-ParagraphXml $cur $xml3
$cur = Insert-ParagraphXml $cur $xml4
$cur = Insert-ParagraphXml $cur $xml5
$cur = Insert-ParagraphXml $cur $xml6
$cur = Insert-ParagraphXml $cur $xml7
$cur = Insert-ParagraphXml $cur $xml8
$cur = Insert-ParagraphXml $cur $xml9
$cur = Insert-ParagraphXml $cur $xml10
$cur = Insert-ParagraphXml $cur $xml11

Write-Output ("final paragraph count=" + $d.Paragraphs.Count)
